$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("W2:W11")
$range.WrapText = $true
$range.Select()
